$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.924259901046753
$ws.Range("B1").Value = 3.028746604919434
$ws.Range("C1").Value = 2.663638114929199
$ws.Range("D1").Value = 2.89653491973877
$ws.Range("E1").Value = 2.836103916168213
